$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for B2:G11 (rows for Q0..Q9), only numeric columns B-G changed.
$data = @(
    @(0.08410299302171197, 1.253677763513526, 6.261918074916323, 2.502382479741321, 2.531284644362528, 42),
    @(-0.01946844920150798, 1.273117676343247, 6.739870713701541, 2.596126097419296, 2.628303442435235, 41),
    @(0.0393246434631094, 1.235849058089852, 6.613468431034368, 2.571666469632944, 2.604123358268155, 40),
    @(-0.04878613897638684, 1.266707606888937, 6.697608253045909, 2.587973773639507, 2.621339043423467, 39),
    @(0.01069592955521734, 1.306387404899978, 6.984373269163457, 2.642796486520189, 2.678249917277403, 38),
    @(-0.03495980545002601, 1.341848499270465, 6.976450903423854, 2.64129720088896, 2.677496044559739, 37),
    @(0.04415342685500446, 1.313189931210337, 7.030782837840379, 2.651562339044733, 2.688802172282798, 36),
    @(-0.01865779601231733, 1.387484845158188, 7.334554756911054, 2.708238312429513, 2.747711521186984, 35),
    @(0.0924138300119391, 1.344817476870875, 7.384053110940232, 2.717361424422639, 2.756630766477241, 34),
    @(-0.01024319886757436, 1.413099074990776, 7.676407440777855, 2.770633039718153, 2.813571923966318, 33)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
